# Applies updated crypto price/volume figures (D/E columns) to sheet1,
# matching the GitHub Actions data refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.141.23'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '1.573.69'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.22'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.796.15'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '1.580.72'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = '27.149.05'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.39'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.31'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '0.0₃0684'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.43'
$ws.Range("E23").Value = '  -3.37%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.63'
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  -3.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.93'
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("D33").Value = '1.394.66'
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("E35").Value = '  +1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.29'
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.939'
$ws.Range("E37").Value = '  -3.20%  '
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  +3.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.81'
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.43'
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.19'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.75'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = '1.709.05'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.59'
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").Value = '0.0₇0989'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0498'
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0952'
$ws.Range("E51").Value = '  -0.62%  '
